# Update the fitting parameters (r_s_star, h_p_star) on the "Parameters"
# sheet, columns J/K of row 2, per the latest Settling-Tool fit results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

$ws.Range("J2").Value = 0.01005
$ws.Range("K2").Value = 0.241
